$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.538.15"
$ws.Range("E2").Value = "  +4.14%  "
$ws.Range("D3").Value = "2.464.73"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'322.83"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'105.12"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("D7").Value = "'0.522"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "'36.12"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'18.29"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("D15").Value = "2.854.07"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "2.513.16"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").Value = "'0.844"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "46.475.78"
$ws.Range("D19").Value = "'12.69"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").Value = "'6.45"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").Value = "'248.75"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  +4.34%  "
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").Value = "'26.10"
$ws.Range("E26").Value = "  +3.55%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'2.31"
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("E30").Value = "  +5.01%  "
$ws.Range("D31").Value = "'49.58"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("D33").Value = "'19.58"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "'5.32"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.0767"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "'2.94"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").Value = "'123.07"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "'20.67"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "1.982.61"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "'2.98"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D48").Value = "'1.79"
$ws.Range("E48").Value = "  +3.99%  "
$ws.Range("E49").Value = "  -4.81%  "
$ws.Range("D50").Value = "'5.30"
$ws.Range("E50").Value = "  +14.52%  "
$ws.Range("D51").Value = "'79.08"
$ws.Range("E51").Value = "  +5.18%  "
